# Insert 3 new rows before row 873 (shifting old 873-918 down to 876-921)
# and populate the new rows with the latest weekly price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertRange = $ws.Range("A873:R875")
$insertRange.Insert(-4121) | Out-Null   # xlShiftDown = -4121

# Copy the date style (s="2", a date/time number format) from the row that
# used to be at 873 (now at 876) down into the 3 new rows, so D873:D875 keep
# the same date formatting as the rest of the column.
$ws.Range("D876").Copy() | Out-Null
$ws.Range("D873:D875").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

function Set-Row($r, $fecha, $calidad, $volumen, $pmin, $pmax, $pprom) {
    $ws.Cells.Item($r, 1).Value = 6
    $ws.Cells.Item($r, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = 100114014
    $ws.Cells.Item($r, 7).Value = "Betarraga"
    $ws.Cells.Item($r, 8).Value = "Sin especificar"
    $ws.Cells.Item($r, 9).Value = $calidad
    $ws.Cells.Item($r, 10).Value = $volumen
    $ws.Cells.Item($r, 11).Value = $pmin
    $ws.Cells.Item($r, 12).Value = $pmax
    $ws.Cells.Item($r, 13).Value = $pprom
    $ws.Cells.Item($r, 14).Value = "`$/unidad"
    $ws.Cells.Item($r, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($r, 16).Value = $pprom
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}

Set-Row 873 44706 "Primera" 50000 110 115 112
Set-Row 874 44706 "Segunda" 44000 85 90 87
Set-Row 875 44706 "Tercera" 15000 65 65 65
